$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking values (e.g. "1.002") as plain
# text in this workbook. Pre-format the specific cells whose price changes
# as Text so Excel does not silently convert them to numbers on assignment.
$priceCells = @(2, 3, 4, 6, 7, 8, 9, 10, 12, 14, 15, 16, 17, 18, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 48, 49, 50, 51)
foreach ($r in $priceCells) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "26.387.66"
$ws.Cells.Item(2, 5).Value = "  -2.36%  "

$ws.Cells.Item(3, 4).Value = "1.776.33"
$ws.Cells.Item(3, 5).Value = "  -0.99%  "

$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.89%  "

$ws.Cells.Item(5, 5).Value = "  -0.57%  "

$ws.Cells.Item(6, 4).Value = "306.49"
$ws.Cells.Item(6, 5).Value = "  -0.31%  "

$ws.Cells.Item(7, 4).Value = "0.4235"
$ws.Cells.Item(7, 5).Value = "  +1.96%  "

$ws.Cells.Item(8, 4).Value = "0.3605"
$ws.Cells.Item(8, 5).Value = "  +2.01%  "

$ws.Cells.Item(9, 4).Value = "0.07153"
$ws.Cells.Item(9, 5).Value = "  +2.13%  "

$ws.Cells.Item(10, 4).Value = "0.8375"
$ws.Cells.Item(10, 5).Value = "  -0.12%  "

$ws.Cells.Item(11, 5).Value = "  +1.94%  "

$ws.Cells.Item(12, 4).Value = "1.763.90"
$ws.Cells.Item(12, 5).Value = "  -6.47%  "

$ws.Cells.Item(13, 5).Value = "  +2.23%  "

$ws.Cells.Item(14, 4).Value = "5.248"
$ws.Cells.Item(14, 5).Value = "  +0.31%  "

$ws.Cells.Item(15, 4).Value = "0.06892"
$ws.Cells.Item(15, 5).Value = "  +1.15%  "

$ws.Cells.Item(16, 4).Value = "1.001"
$ws.Cells.Item(16, 5).Value = "  -0.98%  "

$ws.Cells.Item(17, 4).Value = "79.08"
$ws.Cells.Item(17, 5).Value = "  -0.45%  "

$ws.Cells.Item(18, 4).Value = "0.000008664"
$ws.Cells.Item(18, 5).Value = "  -0.03%  "

$ws.Cells.Item(19, 5).Value = "  -0.79%  "

$ws.Cells.Item(20, 4).Value = "14.90"
$ws.Cells.Item(20, 5).Value = "  -0.50%  "

$ws.Cells.Item(21, 4).Value = "26.396.65"
$ws.Cells.Item(21, 5).Value = "  -3.73%  "

$ws.Cells.Item(22, 4).Value = "5.091"
$ws.Cells.Item(22, 5).Value = "  +1.30%  "

$ws.Cells.Item(23, 4).Value = "10.88"
$ws.Cells.Item(23, 5).Value = "  +2.24%  "

$ws.Cells.Item(24, 4).Value = "1.984.71"
$ws.Cells.Item(24, 5).Value = "  -3.78%  "

$ws.Cells.Item(25, 4).Value = "151.73"
$ws.Cells.Item(25, 5).Value = "  -0.45%  "

$ws.Cells.Item(26, 4).Value = "1.800"
$ws.Cells.Item(26, 5).Value = "  -7.58%  "

$ws.Cells.Item(27, 4).Value = "17.98"
$ws.Cells.Item(27, 5).Value = "  -0.43%  "

$ws.Cells.Item(28, 4).Value = "5.098"
$ws.Cells.Item(28, 5).Value = "  +2.20%  "

$ws.Cells.Item(29, 4).Value = "114.28"
$ws.Cells.Item(29, 5).Value = "  +1.86%  "

$ws.Cells.Item(30, 4).Value = "1.840"
$ws.Cells.Item(30, 5).Value = "  +12.33%  "

$ws.Cells.Item(31, 4).Value = "0.08828"
$ws.Cells.Item(31, 5).Value = "  -0.11%  "

$ws.Cells.Item(32, 4).Value = "0.7283"
$ws.Cells.Item(32, 5).Value = "  +1.78%  "

$ws.Cells.Item(33, 4).Value = "1.124"
$ws.Cells.Item(33, 5).Value = "  +5.14%  "

$ws.Cells.Item(34, 5).Value = "  -0.01%  "

$ws.Cells.Item(35, 4).Value = "1.000"
$ws.Cells.Item(35, 5).Value = "  -0.86%  "

$ws.Cells.Item(36, 4).Value = "2.731"
$ws.Cells.Item(36, 5).Value = "  -5.19%  "

$ws.Cells.Item(37, 4).Value = "1.090"
$ws.Cells.Item(37, 5).Value = "  +1.71%  "

$ws.Cells.Item(38, 4).Value = "0.05110"
$ws.Cells.Item(38, 5).Value = "  +0.84%  "

$ws.Cells.Item(39, 4).Value = "0.01885"
$ws.Cells.Item(39, 5).Value = "  +0.27%  "

$ws.Cells.Item(40, 4).Value = "0.1608"
$ws.Cells.Item(40, 5).Value = "  +0.32%  "

$ws.Cells.Item(41, 4).Value = "0.4912"
$ws.Cells.Item(41, 5).Value = "  +0.15%  "

$ws.Cells.Item(42, 4).Value = "2.600"
$ws.Cells.Item(42, 5).Value = "  -0.99%  "

$ws.Cells.Item(43, 4).Value = "6.329"
$ws.Cells.Item(43, 5).Value = "  +3.25%  "

$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(44, 4).Value = "8.034"
$ws.Cells.Item(44, 5).Value = "  +0.43%  "

$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(45, 4).Value = "10.23"
$ws.Cells.Item(45, 5).Value = "  +0.21%  "

$ws.Cells.Item(46, 4).Value = "104.59"
$ws.Cells.Item(46, 5).Value = "  +1.02%  "

$ws.Cells.Item(47, 5).Value = "  -0.88%  "

$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 4).Value = "1.628"
$ws.Cells.Item(48, 5).Value = "  +3.32%  "

$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).Value = "0.06179"
$ws.Cells.Item(49, 5).Value = "  -2.06%  "

$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).Value = "0.4446"
$ws.Cells.Item(50, 5).Value = "  -1.23%  "

$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).Value = "1.721"
$ws.Cells.Item(51, 5).Value = "  +4.46%  "
